$d = $word.ActiveDocument

# 1) Collapse the "Mein Leben..." sentence (spread across several runs with
#    proofErr spell/grammar markers) into a single plain run of text.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "Mein Leben dreht sich um Lernen, Arbeiten(man muss Geld verdienen) und Schlafen. Eine franzözische Redewendung dafür ist:Metro-Boulot-Dodo.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Mein Leben dreht sich um Lernen, Arbeiten(man muss Geld verdienen) und Schlafen. Eine franzözische Redewendung dafür ist:Metro-Boulot-Dodo.",
    2
)

# 2) Replace the closing line "Mit freundlichen Grüßen" with "MFG"
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute(
    "Mit freundlichen Grüßen", $true, $false, $false, $false, $false, $true, 1, $false,
    "MFG", 2
)
